# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, rows 2-6.
$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 1967
    3 = 592
    4 = 1303
    5 = 6515
    6 = 160
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
